$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 and A3 with new names
$ws.Range("A2").Value = "မန်ကျီးပင်"
$ws.Range("A3").Value = "ရေကျော်"

# Remove A4's content entirely (row 4 element disappears from the sheet)
$ws.Range("A4").ClearContents()

# A5 becomes an empty string value (cell remains present but empty).
# A plain "" assignment is normalized away to a fully blank/absent cell,
# so use a formula that evaluates to the empty string to keep the cell
# present with an empty text value.
$ws.Range("A5").Formula = '=""'
